$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Insert a new row at row 90 (product list is sorted alphabetically; the new
# "SABAO EM BARRA" item sorts just before the existing "SABAO EM PO" row),
# shifting rows 90-108 down to 91-109.
$ws.Rows.Item(90).Insert()

# Populate the new row with the new product data.
$ws.Cells.Item(90, 1).Value = "SABAO EM BARRA - 200G - 200G"
$ws.Cells.Item(90, 2).Value = "UN"
$ws.Cells.Item(90, 3).Value = "S010046"
$ws.Cells.Item(90, 4).Value = 51

# The sheet's hidden filter-database defined name needs to grow to cover the
# newly inserted row.
$wb.Names.Item("Planilha1!_FilterDatabase").RefersTo = "=Planilha1!`$A`$1:`$D`$109"

# The conditional formatting ("notContainsBlanks") range also needs to grow
# to cover the newly inserted row.
$cf = $ws.Range("A2:D108").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A2:D109"))

# Update the active selection to match the saved view.
$ws.Range("E96").Select()
